$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text = "88×13=1144"
$tbl.Cell(1, 2).Range.Text = "97×96=9312"
$tbl.Cell(1, 3).Range.Text = "96×25=2400"
$tbl.Cell(1, 4).Range.Text = "16×25=400"
$tbl.Cell(1, 5).Range.Text = "30×82=2460"
$tbl.Cell(5, 1).Range.Text = "83×85=7055"
$tbl.Cell(5, 2).Range.Text = "27×75=2025"
$tbl.Cell(5, 3).Range.Text = "58×32=1856"
$tbl.Cell(5, 4).Range.Text = "88×48=4224"
$tbl.Cell(5, 5).Range.Text = "91×18=1638"
$tbl.Cell(10, 1).Range.Text = "21×25=525"
$tbl.Cell(10, 2).Range.Text = "46×97=4462"
$tbl.Cell(10, 3).Range.Text = "28×89=2492"
$tbl.Cell(10, 4).Range.Text = "82×70=5740"
$tbl.Cell(10, 5).Range.Text = "92×94=8648"
$tbl.Cell(15, 1).Range.Text = "52×65=3380"
$tbl.Cell(15, 2).Range.Text = "94×50=4700"
$tbl.Cell(15, 3).Range.Text = "42×94=3948"
$tbl.Cell(15, 4).Range.Text = "26×51=1326"
$tbl.Cell(15, 5).Range.Text = "53×92=4876"
$tbl.Cell(20, 1).Range.Text = "23×46=1058"
$tbl.Cell(20, 2).Range.Text = "56×18=1008"
$tbl.Cell(20, 3).Range.Text = "62×38=2356"
$tbl.Cell(20, 4).Range.Text = "23×16=368"
$tbl.Cell(20, 5).Range.Text = "99×75=7425"
